$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.159.46'
$ws.Range('E2').Value = '  +1.23%  '
$ws.Range('D3').Value = '3.510.07'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '601.78'
$ws.Range('E5').Value = '  -0.90%  '
$ws.Range('D6').Value = '174.32'
$ws.Range('E6').Value = '  +1.86%  '
$ws.Range('D7').Value = '0.610'
$ws.Range('E7').Value = '  -0.59%  '
$ws.Range('D8').Value = '3.504.45'
$ws.Range('E8').Value = '  -0.32%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  -1.42%  '
$ws.Range('D11').Value = '7.23'
$ws.Range('E11').Value = '  +7.98%  '
$ws.Range('D12').Value = '0.581'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('D13').Value = '46.20'
$ws.Range('E13').Value = '  -2.30%  '
$ws.Range('D14').Value = '0.0000275'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').Value = '4.068.76'
$ws.Range('E15').Value = '  -0.36%  '
$ws.Range('D16').Value = '612.99'
$ws.Range('E16').Value = '  -0.65%  '
$ws.Range('D17').Value = '8.28'
$ws.Range('E17').Value = '  -0.93%  '
$ws.Range('D18').Value = '3.504.49'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').Value = '70.194.70'
$ws.Range('E19').Value = '  +1.12%  '
$ws.Range('D20').Value = '0.120'
$ws.Range('E20').Value = '  +0.53%  '
$ws.Range('D21').Value = '17.24'
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('D22').Value = '0.874'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('D23').Value = '9.08'
$ws.Range('E23').Value = '  -19.15%  '
$ws.Range('D24').Value = '15.53'
$ws.Range('E24').Value = '  -1.79%  '
$ws.Range('D25').Value = '96.04'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').Value = '3.70'
$ws.Range('E26').Value = '  -4.53%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').Value = '2.56'
$ws.Range('E28').Value = '  -2.65%  '
$ws.Range('D29').Value = '34.09'
$ws.Range('E29').Value = '  +2.28%  '
$ws.Range('D30').Value = '8.96'
$ws.Range('E30').Value = '  -3.14%  '
$ws.Range('D31').Value = '8.12'
$ws.Range('E31').Value = '  -4.43%  '
$ws.Range('E32').Value = '  -5.47%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '639.79'
$ws.Range('E33').Value = '  +12.08%  '
$ws.Range('B34').Value = 'Mantle'
$ws.Range('C34').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D34').Value = '1.28'
$ws.Range('E34').Value = '  -4.08%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '6.85'
$ws.Range('E35').Value = '  -0.92%  '
$ws.Range('D36').Value = '3.57'
$ws.Range('E36').Value = '  +0.44%  '
$ws.Range('D37').Value = '0.0993'
$ws.Range('E37').Value = '  -2.49%  '
$ws.Range('D38').Value = '10.69'
$ws.Range('E38').Value = '  -0.84%  '
$ws.Range('D39').Value = '0.0472'
$ws.Range('E39').Value = '  +6.08%  '
$ws.Range('B40').Value = 'OKB'
$ws.Range('C40').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D40').Value = '56.56'
$ws.Range('E40').Value = '  -1.15%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('D42').Value = '0.142'
$ws.Range('E42').Value = '  +1.76%  '
$ws.Range('D43').Value = '0.0₃0741'
$ws.Range('E43').Value = '  +4.62%  '
$ws.Range('D44').Value = '3.348.63'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('D45').Value = '0.310'
$ws.Range('E45').Value = '  -5.08%  '
$ws.Range('D46').Value = '2.91'
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').Value = '32.10'
$ws.Range('E47').Value = '  -2.51%  '
$ws.Range('D48').Value = '2.55'
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').Value = '134.36'
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('E51').Value = '  +0.00%  '
